$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "-"

$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "[-, -, 'MCT-3A-Eletrohidráulica', -]"
$ws.Range("F3").Value = "[-, -, 'MEC-3B-Coman. Hidraulicos', -]"

$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "[-, -, 'MCT-3A-Eletrohidráulica', -]"
$ws.Range("F4").Value = "[-, -, 'MEC-3B-Coman. Hidraulicos', -]"

$ws.Range("E6").Value = "[-, -, 'MCT-3A-Eletrohidráulica', -]"
$ws.Range("F6").Value = "[-, -, 'MEC-3B-Coman. Hidraulicos', -]"

$ws.Range("E7").Value = "[-, -, 'MCT-3A-Eletrohidráulica', -]"
$ws.Range("F7").Value = "[-, -, 'MEC-3B-Coman. Hidraulicos', -]"

$ws.Range("F8").Value = "-"

$ws.Range("B10").Value = "-"

$ws.Range("B11").Value = "[-, 'MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP']"
$ws.Range("F11").Value = "MEC-2A-Máquinas Térmicas e de Fluxo"

$ws.Range("B12").Value = "[-, 'MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP']"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "MEC-2A-Máquinas Térmicas e de Fluxo"

$ws.Range("B14").Value = "[-, 'MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP']"
$ws.Range("B15").Value = "[-, 'MEC-3A-Comandos Eletricos', -, 'MEC-3A-Cont.Lóg.Prog CLP']"

$ws.Range("B16").Value = "-"
